# Journal de bord update: document row 96 and 97 with new journal entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 96 -----------------------------------------------------------
$ws.Range("A96").Value = "Documentation"
$ws.Range("A96").WrapText = $true

$ws.Range("B96").Value = "Documentation du code aux endroits où j'estimais qu'il y avait un manque de clarté"
$ws.Range("B96").WrapText = $true

$ws.Range("C96").Value = 1

$ws.Range("D96").Value = 43558

$ws.Rows(96).RowHeight = 45

# --- Row 97 -----------------------------------------------------------
$ws.Range("A97").Value = "Documentation"
$ws.Range("A97").WrapText = $true

$ws.Range("B97").Value = "Vérification de la présence de tous les use cases , clarification de certains use cases / scénarios. Ajout de scénarios manquant + tests sur l'API"
$ws.Range("B97").WrapText = $true

$ws.Range("C97").Value = 2

$ws.Range("D97").Value = 43558

$ws.Rows(97).RowHeight = 75

# --- Selection / view ---------------------------------------------------
$ws.Range("C96").Select()

Write-Host "Journal de bord rows 96-97 updated"
